$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: effort on B33 changed from 2.75 to 4.75, and the "Additional Effort"
# entry in C33 (2) is removed entirely.
$ws.Range("B33").Value = 4.75
$ws.Range("C33").Value = $null

# New last row (58): a new log entry dated 13/12/2012, 0.5h effort and a
# descriptive comment in column D.
$ws.Range("A58").Value = "12/13/2012"
$ws.Range("B58").Value = 0.5
$ws.Range("D58").Value = "Documentation of code slightly improved"

# Reflect the new selection left behind by the edit (C33 was last touched).
$ws.Range("C33").Select()
